$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember the width of the column immediately to the left (M) so the
# newly-inserted column (N) can inherit it, matching Excel's native
# "insert column" behaviour of copying the left neighbour's formatting.
$leftColumnWidth = $ws.Columns("M").ColumnWidth

# Insert a new blank column before column N ("Late"), shifting the
# "Late"/"Outstanding" (and the second "heading" helper) columns one to
# the right.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $leftColumnWidth

# Activate the "Repayment schedule" sheet/tab and select cell S3, matching
# the new active tab/selection recorded in the workbook.
$ws.Activate()
$ws.Range("S3").Select()
